$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Build the two new border-only cell styles on sheet1 (quality_comparison) ---
# C1: top + bottom thin borders (matches existing border index 4)
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# D1: right + top + bottom thin borders (matches existing border index 5)
$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$d1.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$d1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# --- Reuse the exact same formats (copy/paste format) on sheet2 (computational_comparison) ---
# This avoids creating duplicate/unused style entries in cellXfs.
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$c1.Copy()
$ws2.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$d1.Copy()
$ws2.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Text replacements: "fedcore" -> "approach" ---
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Remove the stray empty inline-string cell G5 on sheet2 ---
$ws2.Range("G5").ClearContents()

Write-Host "edit applied"
